$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the sheet/tab name to reflect the new "through" date.
$ws.Name = "Through 2022-08-13"

# Update the header label in B1 (shared string) to match.
$ws.Range("B1").Value = "August 2022 (through August 13)"

# Cell value updates / additions (row => neighborhood, column => month-year)
$ws.Range("B2").Value = 11
$ws.Range("R3").Value = 2
$ws.Range("B5").Value = 4
$ws.Range("J5").Value = 9
$ws.Range("B6").Value = 6
$ws.Range("B7").Value = 5
$ws.Range("B8").Value = 1
$ws.Range("AP8").Value = 4
$ws.Range("BF8").Value = 1
$ws.Range("AX12").Value = 1
$ws.Range("R14").Value = 1
$ws.Range("AX22").Value = 1
$ws.Range("R23").Value = 1
$ws.Range("R28").Value = 2
$ws.Range("R32").Value = 1
$ws.Range("R35").Value = 1
$ws.Range("B42").Value = 2
$ws.Range("J45").Value = 2
$ws.Range("AP50").Value = 1
$ws.Range("B54").Value = 2
$ws.Range("B60").Value = 1
$ws.Range("J60").Value = 1
$ws.Range("R64").Value = 3
$ws.Range("AX64").Value = 1
